# "Changed format to Jupyter notebook"
#
# Adds two new paragraphs after the last "Took <Dataset> from Kaggle
# dataset" note (the Cleveland Heart Disease Dataset paragraph) and before
# the document's trailing (pre-existing) empty paragraphs:
#   1. an empty bold paragraph (a spacer), then
#   2. a bold paragraph reading:
#      "Create conda environment to install sepereate dependecies"

$d = $word.ActiveDocument

# Locate the "Took Cleveland Heart Disease Dataset from Kaggle dataset"
# paragraph - the last of the three dataset notes, and the point after
# which the new paragraphs belong.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Cleveland Heart Disease Dataset*") {
        $targetIndex = $i
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the 'Cleveland Heart Disease Dataset' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Collapse to the end of that paragraph (i.e. right after its paragraph
# mark) and insert a new, empty paragraph there. This becomes the first of
# the two new paragraphs (the blank spacer line).
$insertPoint = $target.Range.Duplicate
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()

# The paragraph immediately after the spacer is where the second new
# paragraph's text goes; split again right after the spacer's mark and
# type the note text there.
$spacer = $d.Paragraphs.Item($targetIndex + 1)
$textPoint = $spacer.Range.Duplicate
$textPoint.Collapse(0)
$textPoint.InsertParagraphAfter()
$textPoint.InsertAfter("Create conda environment to install sepereate dependecies")

# Make sure the new note paragraph is bold (matching the rest of the
# document's notes), including the complex-script bold flag.
$notePara = $d.Paragraphs.Item($targetIndex + 2)
$notePara.Range.Font.Bold = -1
$notePara.Range.Font.BoldBi = -1

Write-Output "Document now has $($d.Paragraphs.Count) paragraphs."
